# Record deleted: 11364560
# Row A1319 ("11364560", MAKS, HAVVA NİLGÜN KIYMAÇ / SEVAL ÇELİK) is removed
# from the "Kayitlar" sheet, and the same record (mirrored per-district on
# the "Merkez İlçe" sheet) is removed from there as well. Deleting the
# entire row shifts all subsequent rows up by one, matching the diff.

$wb = $excel.ActiveWorkbook

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows.Item(1319).Delete()

$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows.Item(780).Delete()
